$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "Volume" header in F1 (and clear the column)
$ws.Range("F1").ClearContents()

# Add new data rows for sequential experiments
$ws.Range("A2").Value = 1
$ws.Range("B2").Value = 50
$ws.Range("C2").Value = 10

$ws.Range("A3").Value = 2
$ws.Range("B3").Value = 10
$ws.Range("C3").Value = 50

$ws.Range("A4").Value = 3
$ws.Range("C4").Value = 10
$ws.Range("D4").Value = 50

# Update the selection to match the target state
$ws.Range("D10").Select()
